$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Data for subjects 21-26 (rows 23-28), columns C (Q1) through Q (Q15)
$data = @{
    23 = @(5,4,4,4,5,2,5,5,4,2,5,5,5,2,6)
    24 = @(4,2,4,1,2,2,4,3,3,4,5,4,5,5,5)
    25 = @(4,4,5,4,5,5,3,4,2,4,4,4,4,5,4)
    26 = @(3,5,4,3,3,5,4,5,2,4,5,5,3,5,6)
    27 = @(4,3,2,3,4,3,4,5,4,3,3,3,5,2,6)
    28 = @(4,2,3,1,1,1,4,1,3,2,3,3,3,3,4)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 3 + $i  # Column C = 3
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}

$ws.Range("Q29").Select()
